# "Time to start working" — duplicate the current Methods section as the
# new Methods section (leaving the "gathered from" sentence unfinished),
# and turn the old Methods section's heading into Results, replacing its
# old body paragraphs with a new Discussion heading.

$d = $word.ActiveDocument

# The 3rd paragraph is the "Introduction" body paragraph; the 4th is the
# existing "Methods" heading we are about to duplicate.
$introBodyPara = $d.Paragraphs.Item(3)

# --------------------------------------------------------------------
# 1) Insert a brand-new "Methods" section right after the Introduction
#    paragraph (i.e. before the existing Methods heading), duplicating
#    the text currently in the Methods section below it.
# --------------------------------------------------------------------

# Heading: "Methods" (bold)
$introBodyPara.Range.InsertParagraphAfter()
$newMethodsHeading = $d.Paragraphs.Item(4)
$newMethodsHeading.Range.Text = "Methods"
$newMethodsHeading.Range.Font.Bold = $true

# Paragraph: "To be used as parameters..." left as an unfinished sentence
$newMethodsHeading.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item(5)
$p1.Range.Font.Bold = $false
$p1.Range.Text = "To be used as parameters, datasets on cost of living, population density, and light pollution (current list, more will be added soon) were gathered from … (not a finished sentence)"

# Paragraph: "Since datasets for this study..."
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item(6)
$p2.Range.Text = "Since datasets for this study aren’t that big – in most cases couple hundreds of countries, they were manually cleaned from any unnecessary data and checked for having the same number of countries and that the countries are the same."

# Paragraph: "To make my system work..."
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item(7)
$p3.Range.Text = "To make my system work, necessary parameters will be chosen by system’s user in the beginning and acceptable value ranges will be chosen to create a custom(personalized) decision tree to help with categorizing our data into “good” and “bad”. "

# Paragraph: "To further refine our results..."
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(8)
$p4.Range.Text = "To further refine our results, two separate systems will be used – priority list and k-means clustering algorithm. When user will select the parameters for their research, the user will be prompted to assign a priority to each of the parameters they chose. It will allow to favor certain parameters compared to others when time comes to calculate countries “priority points”, which will calculate a value for each of the countries based on priorities set for different parameters, thus allowing to rank countries from “best” to “worst” in their respective categories (“good” or “bad”). "

# Paragraph: "To categorize countries in a vaguer way..."
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(9)
$p5.Range.Text = "To categorize countries in a vaguer way, excluding the priorities of the parameters, k-means clustering will be used to categorize countries in three distinct groups: “best fit”, “medium fit”, “worst fit”. This will create more of the suggestive categories for people to observe. For example: let’s take our theoretical “good” group of countries and run it through k-means algorithm; this will yield three subgroups of “best fit”, “medium fit”, and “worst fit” countries from our master-group “good countries”."

# Trailing blank paragraph, matching the blank line that used to sit
# right before the next (now pushed-down) section heading.
$p5.Range.InsertParagraphAfter()

# --------------------------------------------------------------------
# 2) The original "Methods" heading (now pushed further down the
#    document, at index 11) becomes "Results". Its old body paragraphs
#    (indices 12-16) are removed and replaced by a "Discussion" heading.
# --------------------------------------------------------------------

$oldMethodsHeading = $d.Paragraphs.Item(11)
$oldMethodsHeading.Range.Text = "Results"
$oldMethodsHeading.Range.Font.Bold = $true

$oldBodyStart = $d.Paragraphs.Item(12).Range.Start
$oldBodyEnd = $d.Paragraphs.Item(16).Range.End
$oldBodyRange = $d.Range($oldBodyStart, $oldBodyEnd)
$oldBodyRange.Delete()

$oldMethodsHeading.Range.InsertParagraphAfter()
$discussionPara = $d.Paragraphs.Item(12)
$discussionPara.Range.Text = "Discussion"
$discussionPara.Range.Font.Bold = $true
